## Add files via upload
## ---------------------
## The single "Sheet1" (Name/DormNum/StudentID, 1 data row) becomes three
## grade-level roster sheets: "1stGrade", "2ndGrade", "3rdGrade" - each with
## the same DormNum/Name/StudentID header and two student rows.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and append two more after it ---------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "1stGrade"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2ndGrade"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "3rdGrade"

# --- 1stGrade: keep existing header, add the two student rows -------------
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "김XX"
$ws1.Range("C2").Value = 11001

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "박XX"
$ws1.Range("C3").Value = 10320

# --- 2ndGrade: new header + two student rows -------------------------------
$ws2.Range("A1").Value = "DormNum"
$ws2.Range("B1").Value = "Name"
$ws2.Range("C1").Value = "StudentID"

$ws2.Range("A2").Value = 2
$ws2.Range("B2").Value = "구OO"
$ws2.Range("C2").Value = 20912

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "이OO"
$ws2.Range("C3").Value = 20311

# --- 3rdGrade: new header + two student rows -------------------------------
$ws3.Range("A1").Value = "DormNum"
$ws3.Range("B1").Value = "Name"
$ws3.Range("C1").Value = "StudentID"

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "박YY"
$ws3.Range("C2").Value = 30910

$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "김YY"
$ws3.Range("C3").Value = 30123

# --- Restore per-sheet selections as recorded in the saved workbook -------
$ws1.Range("A2:C3").Select()
$ws2.Range("C3").Select()
$ws3.Range("B3").Select()

# --- 3rdGrade is the tab that was active/visible when the file was saved --
$ws3.Activate()
